$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the steel (S) description line to remove the "RME/" segment.
$cell = $ws.Range("B2")
$cell.Value = "6% CR/LFM+CDM/H:1`n31% CR/LWAL+CDH/H:1`n11% CR+PC/LFM+CDL/H:1`n9% CR/LWAL+CDM/HBET:3-5`n12% S/LFM+CDH/H:1`n1% S/LFM+CDH/HBET:3-5`n14% S/LFBR+CDH/H:1`n14% W/LWAL+CDM/H:1`n1% MUR/LWAL+CDN/H:1`n1% MR/LWAL+CDL/H:1"

# Wrap text so the multi-line description is fully visible.
$cell.WrapText = $true

# Expand the row to the maximum height to show all wrapped lines.
$ws.Rows.Item(2).RowHeight = 409.6

# Match the saved selection state.
$ws.Range("C2:C13").Select()
